$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 / Row 3: account holder name + card number ---
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long digit string that must stay literal TEXT (it was already
# stored as text before the edit) without disturbing its existing style.
# Writing a numeric-looking string straight into Value coerces it to a
# Number, so stage it on a scratch cell formatted as Text, copy just the
# VALUE over (PasteSpecial values-only leaves B3's own formatting intact),
# then wipe the scratch cell completely so no residue is left behind.
$scratch = $ws.Range("G3")
$scratch.NumberFormat = "@"
$scratch.Value = "2570314725427075"
$scratch.Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null

$ws.Range("C3").Value = "Mohaupt"

# --- Row 5: opening balance date ---
$ws.Range("D5").Value = "KONTOSTAND AM 22.04.2024"

# --- Row 6: transaction 1 ---
$ws.Range("B6").Value = "24.04."
$ws.Range("C6").Value = "25.04."
$ws.Range("D6").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 32144291"
$ws.Range("E6").Value = "84,95-"

# --- Row 7: transaction 2 ---
$ws.Range("B7").Value = "26.04."
$ws.Range("C7").Value = "27.04."
$ws.Range("D7").Value = "BEITRAG Allianz SE K-33763106"
$ws.Range("E7").Value = "57,46-"

# --- Row 8: transaction 3 ---
$ws.Range("B8").Value = "29.04."
$ws.Range("C8").Value = "30.04."
$ws.Range("D8").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E8").Value = "24,93-"

# --- Row 9: transaction 4 (was blank, now a new row of data) ---
$ws.Range("B9").Value = "03.05."
$ws.Range("C9").Value = "04.05."
$ws.Range("D9").Value = "MCDONALDS Viersen"
# E9 needs the right-aligned amount style used by the other amount cells
# (matches E8's formatting) instead of the blank placeholder style.
$ws.Range("E8").Copy() | Out-Null
$ws.Range("E9").PasteSpecial(-4122) | Out-Null
$ws.Range("E9").Value = "38,15-"

# --- Row 12: closing balance ---
$ws.Range("D12").Value = "KONTOSTAND AM 08.05.2024"
$ws.Range("E12").Value = "205,49-"

# --- Row 13: next statement date ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 17.05.2024"
